$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number that was bumped by
# one day (45188 -> 45189, i.e. 2023-09-19 -> 2023-09-20) for every data
# row (rows 2 through 301).
$ws.Range("C2:C301").Value = 45189
